$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 347 (shifts existing rows 347-384 down to 348-385)
$ws.Rows.Item(347).Insert()

# Populate the new row 347 with the new weekly data point
$ws.Cells.Item(347, 1).Value2 = 7
$ws.Cells.Item(347, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(347, 3).Value2 = "Ñuble"
$ws.Cells.Item(347, 4).Value2 = 45212
$ws.Cells.Item(347, 5).Value2 = 16
$ws.Cells.Item(347, 6).Value2 = 100112043
$ws.Cells.Item(347, 7).Value2 = "Pepino ensalada"
$ws.Cells.Item(347, 8).Value2 = "Sin especificar"
$ws.Cells.Item(347, 9).Value2 = "Primera"
$ws.Cells.Item(347, 10).Value2 = 60
$ws.Cells.Item(347, 11).Value2 = 15000
$ws.Cells.Item(347, 12).Value2 = 16000
$ws.Cells.Item(347, 13).Value2 = 15500
$ws.Cells.Item(347, 14).Value2 = "$/caja 60 unidades"
$ws.Cells.Item(347, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(347, 16).Value2 = 258
$ws.Cells.Item(347, 17).Value2 = 60
$ws.Cells.Item(347, 18).Value2 = "Hortaliza"
